$d = $word.ActiveDocument

$d.Content.Find.Execute("41×67=", $true, $false, $false, $false, $false, $true, 1, $false, "34×27=", 2) | Out-Null
$d.Content.Find.Execute("38×62=", $true, $false, $false, $false, $false, $true, 1, $false, "85×54=", 2) | Out-Null
$d.Content.Find.Execute("71×70=", $true, $false, $false, $false, $false, $true, 1, $false, "35×72=", 2) | Out-Null
$d.Content.Find.Execute("81×83=", $true, $false, $false, $false, $false, $true, 1, $false, "12×11=", 2) | Out-Null
$d.Content.Find.Execute("71×34=", $true, $false, $false, $false, $false, $true, 1, $false, "37×33=", 2) | Out-Null
$d.Content.Find.Execute("16×89=", $true, $false, $false, $false, $false, $true, 1, $false, "27×17=", 2) | Out-Null
$d.Content.Find.Execute("51×73=", $true, $false, $false, $false, $false, $true, 1, $false, "31×73=", 2) | Out-Null
$d.Content.Find.Execute("23×65=", $true, $false, $false, $false, $false, $true, 1, $false, "31×30=", 2) | Out-Null
$d.Content.Find.Execute("94×31=", $true, $false, $false, $false, $false, $true, 1, $false, "88×11=", 2) | Out-Null
$d.Content.Find.Execute("42×53=", $true, $false, $false, $false, $false, $true, 1, $false, "89×35=", 2) | Out-Null
$d.Content.Find.Execute("37×76=", $true, $false, $false, $false, $false, $true, 1, $false, "69×96=", 2) | Out-Null
$d.Content.Find.Execute("74×28=", $true, $false, $false, $false, $false, $true, 1, $false, "90×42=", 2) | Out-Null
$d.Content.Find.Execute("23×70=", $true, $false, $false, $false, $false, $true, 1, $false, "12×40=", 2) | Out-Null
$d.Content.Find.Execute("66×95=", $true, $false, $false, $false, $false, $true, 1, $false, "95×32=", 2) | Out-Null
$d.Content.Find.Execute("89×95=", $true, $false, $false, $false, $false, $true, 1, $false, "65×28=", 2) | Out-Null
$d.Content.Find.Execute("61×63=", $true, $false, $false, $false, $false, $true, 1, $false, "14×25=", 2) | Out-Null
$d.Content.Find.Execute("44×44=", $true, $false, $false, $false, $false, $true, 1, $false, "63×92=", 2) | Out-Null
$d.Content.Find.Execute("64×67=", $true, $false, $false, $false, $false, $true, 1, $false, "57×27=", 2) | Out-Null
$d.Content.Find.Execute("81×68=", $true, $false, $false, $false, $false, $true, 1, $false, "81×40=", 2) | Out-Null
$d.Content.Find.Execute("68×49=", $true, $false, $false, $false, $false, $true, 1, $false, "97×51=", 2) | Out-Null
$d.Content.Find.Execute("18×89=", $true, $false, $false, $false, $false, $true, 1, $false, "16×33=", 2) | Out-Null
$d.Content.Find.Execute("27×82=", $true, $false, $false, $false, $false, $true, 1, $false, "47×59=", 2) | Out-Null
$d.Content.Find.Execute("11×91=", $true, $false, $false, $false, $false, $true, 1, $false, "37×19=", 2) | Out-Null
$d.Content.Find.Execute("69×44=", $true, $false, $false, $false, $false, $true, 1, $false, "74×91=", 2) | Out-Null
$d.Content.Find.Execute("71×33=", $true, $false, $false, $false, $false, $true, 1, $false, "13×75=", 2) | Out-Null
